# Update column G ("K") values on the active sheet to reflect the
# regenerated save_data (K computed instead of the old Strike# values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 0
    4  = 0
    6  = 3
    7  = 2
    8  = 2
    9  = 0
    10 = 2
    11 = 1
    12 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
